# feat: add 2022-Q3 data
#
# - Inserts a new worksheet "2022-Q3" (positioned right before "2022-Q2"),
#   filled with the quarterly fund-holdings table.
# - Updates the "总计" (totals) summary sheet: a new row for 2022-Q3 is
#   inserted above the existing 2022-Q2 / 2021-Q3 rows, and the running
#   index column is renumbered.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet, positioned before "2022-Q2"
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "2022-Q3"

# Re-fetch sheets by name (avoid stale references around Worksheets.Add)
$src = $wb.Worksheets.Item("2022-Q2")
$dst = $wb.Worksheets.Item("2022-Q3")

# Clone the header-row / index-column formatting from the 2022-Q2 sheet so
# the new sheet matches the existing look (bold + bordered cells).
$src.Range("B1:H1").Copy()
$dst.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$src.Range("A2").Copy()
$dst.Range("A2:A6").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $dst.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Data rows: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0, "012262", "华宝可持续发展混合A", "8.44", "84.59", "3.00", "0.2532", 10),
    @(1, "012263", "华宝可持续发展混合C", "4.74", "84.59", "3.00", "0.1422", 10),
    @(2, "002210", "创金合信量化多因子股票A", "2.39", "91.71", "1.28", "0.0306", 4),
    @(3, "004284", "华宝新优选一年定期开放灵活配置混合", "0.53", "75.29", "2.89", "0.0153", 9),
    @(4, "003865", "创金合信量化多因子股票C", "0.75", "91.71", "1.28", "0.0096", 4)
)

$r = 2
foreach ($row in $rows) {
    $dst.Cells.Item($r, 1).Value = $row[0]

    $dst.Cells.Item($r, 2).NumberFormat = "@"
    $dst.Cells.Item($r, 2).Value = $row[1]

    $dst.Cells.Item($r, 3).NumberFormat = "@"
    $dst.Cells.Item($r, 3).Value = $row[2]

    $dst.Cells.Item($r, 4).NumberFormat = "@"
    $dst.Cells.Item($r, 4).Value = $row[3]

    $dst.Cells.Item($r, 5).NumberFormat = "@"
    $dst.Cells.Item($r, 5).Value = $row[4]

    $dst.Cells.Item($r, 6).NumberFormat = "@"
    $dst.Cells.Item($r, 6).Value = $row[5]

    $dst.Cells.Item($r, 7).NumberFormat = "@"
    $dst.Cells.Item($r, 7).Value = $row[6]

    $dst.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q3 row above the
#    existing 2022-Q2 / 2021-Q3 rows, shifting them down by one and
#    renumbering the index column (A).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 4 is new; clone formatting from the existing A2 (index-column) cell.
$total.Cells.Item(2, 1).Copy()
$total.Cells.Item(4, 1).PasteSpecial(-4122)   # xlPasteFormats

# Shift 2021-Q3 from row 3 -> row 4
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q3"
$total.Cells.Item(4, 3).Value = 9
$total.Cells.Item(4, 4).Value = 7.07

# Shift 2022-Q2 from row 2 -> row 3
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 1
$total.Cells.Item(3, 4).Value = 0.08

# New 2022-Q3 summary row -> row 2
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.45
